$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the column-A label style (bold, bordered, centered) to the new row 11
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rewrite rows 3 and 5-11 with the updated / shifted data
# Row 3: even_MAG-GUT47840.fa
$ws.Cells.Item(3,1).Value = "even_MAG-GUT47840.fa"
$ws.Cells.Item(3,2).Value = 0.002128955069187773
$ws.Cells.Item(3,3).Value = [double]"2.884368018219164e-06"
$ws.Cells.Item(3,4).Value = 0.9926253869066943
$ws.Cells.Item(3,5).Value = 0.0002318820314145452
$ws.Cells.Item(3,6).Value = [double]"2.22028685969334e-14"
$ws.Cells.Item(3,7).Value = 0.00283471577611792
$ws.Cells.Item(3,8).Value = [double]"1.136072549169427e-05"
$ws.Cells.Item(3,9).Value = [double]"2.22028685969334e-14"
$ws.Cells.Item(3,10).Value = [double]"2.22028685969334e-14"
$ws.Cells.Item(3,11).Value = [double]"1.795538279463149e-06"
$ws.Cells.Item(3,12).Value = 0.002113297343875105
$ws.Cells.Item(3,13).Value = [double]"3.879748315285884e-05"
$ws.Cells.Item(3,14).Value = [double]"1.092475770173163e-05"
$ws.Cells.Item(3,15).Value = 0.9926253869066943
$ws.Cells.Item(3,16).Value = "g__Fenollaria"
$ws.Cells.Item(3,17).Value = "g__Fenollaria"

# Row 5: even_MAG-GUT7064.fa
$ws.Cells.Item(5,1).Value = "even_MAG-GUT7064.fa"
$ws.Cells.Item(5,2).Value = 0.7312153239712396
$ws.Cells.Item(5,3).Value = 0.04371230977155079
$ws.Cells.Item(5,4).Value = [double]"4.403331657834807e-05"
$ws.Cells.Item(5,5).Value = 0.07037544759369332
$ws.Cells.Item(5,6).Value = [double]"6.276127602520476e-14"
$ws.Cells.Item(5,7).Value = 0.02003777797342363
$ws.Cells.Item(5,8).Value = 0.0002906417970133319
$ws.Cells.Item(5,9).Value = [double]"6.276127602520472e-14"
$ws.Cells.Item(5,10).Value = [double]"6.276127602520472e-14"
$ws.Cells.Item(5,11).Value = 0.02670007922926321
$ws.Cells.Item(5,12).Value = 0.07058264724870372
$ws.Cells.Item(5,13).Value = 0.03112526995981559
$ws.Cells.Item(5,14).Value = 0.005916469138530498
$ws.Cells.Item(5,15).Value = 0.7312153239712396
$ws.Cells.Item(5,16).Value = "g__Anaerococcus"
$ws.Cells.Item(5,17).Value = "g__Anaerococcus(reject)"

# Row 6: even_MAG-GUT7291.fa
$ws.Cells.Item(6,1).Value = "even_MAG-GUT7291.fa"
$ws.Cells.Item(6,2).Value = 0.8424465234678435
$ws.Cells.Item(6,3).Value = 0.0002259659329396396
$ws.Cells.Item(6,4).Value = [double]"1.811583084658021e-05"
$ws.Cells.Item(6,5).Value = 0.03380488686197823
$ws.Cells.Item(6,6).Value = [double]"3.76558110959485e-14"
$ws.Cells.Item(6,7).Value = 0.01276822636079955
$ws.Cells.Item(6,8).Value = 0.00080575934966557
$ws.Cells.Item(6,9).Value = [double]"3.76558110959485e-14"
$ws.Cells.Item(6,10).Value = [double]"3.76558110959485e-14"
$ws.Cells.Item(6,11).Value = 0.07351893887462946
$ws.Cells.Item(6,12).Value = 0.02278223068749314
$ws.Cells.Item(6,13).Value = 0.003734819065349844
$ws.Cells.Item(6,14).Value = 0.00989453356834148
$ws.Cells.Item(6,15).Value = 0.8424465234678435
$ws.Cells.Item(6,16).Value = "g__Anaerococcus"
$ws.Cells.Item(6,17).Value = "g__Anaerococcus"

# Row 7: even_MAG-GUT91256.fa
$ws.Cells.Item(7,1).Value = "even_MAG-GUT91256.fa"
$ws.Cells.Item(7,2).Value = 0.0003964139812613548
$ws.Cells.Item(7,3).Value = [double]"7.630037483268974e-06"
$ws.Cells.Item(7,4).Value = [double]"8.904527650208722e-07"
$ws.Cells.Item(7,5).Value = 0.001595620383121524
$ws.Cells.Item(7,6).Value = [double]"2.220295376730079e-14"
$ws.Cells.Item(7,7).Value = 0.004205554150421521
$ws.Cells.Item(7,8).Value = 0.0003175017496848924
$ws.Cells.Item(7,9).Value = [double]"2.220295376730079e-14"
$ws.Cells.Item(7,10).Value = [double]"2.220295376730079e-14"
$ws.Cells.Item(7,11).Value = 0.008670147407215806
$ws.Cells.Item(7,12).Value = 0.9812398988648945
$ws.Cells.Item(7,13).Value = [double]"3.196457623362376e-05"
$ws.Cells.Item(7,14).Value = 0.003534378396851786
$ws.Cells.Item(7,15).Value = 0.9812398988648945
$ws.Cells.Item(7,16).Value = "g__Peptoniphilus_A"
$ws.Cells.Item(7,17).Value = "g__Peptoniphilus_A"

# Row 8: even_MAG-GUT91291.fa
$ws.Cells.Item(8,1).Value = "even_MAG-GUT91291.fa"
$ws.Cells.Item(8,2).Value = 0.0003835460830238558
$ws.Cells.Item(8,3).Value = [double]"7.448656938903593e-05"
$ws.Cells.Item(8,4).Value = 0.0003324225578362673
$ws.Cells.Item(8,5).Value = 0.002052420769840596
$ws.Cells.Item(8,6).Value = [double]"2.220057606210968e-14"
$ws.Cells.Item(8,7).Value = 0.006851916614838984
$ws.Cells.Item(8,8).Value = [double]"4.602408547077245e-05"
$ws.Cells.Item(8,9).Value = [double]"2.220057606210968e-14"
$ws.Cells.Item(8,10).Value = [double]"2.220057606210968e-14"
$ws.Cells.Item(8,11).Value = 0.0004110835217260299
$ws.Cells.Item(8,12).Value = 0.9836877822335476
$ws.Cells.Item(8,13).Value = 0.002358505238309923
$ws.Cells.Item(8,14).Value = 0.003801812325950353
$ws.Cells.Item(8,15).Value = 0.9836877822335476
$ws.Cells.Item(8,16).Value = "g__Peptoniphilus_A"
$ws.Cells.Item(8,17).Value = "g__Peptoniphilus_A"

# Row 9: even_MAG-GUT91328.fa
$ws.Cells.Item(9,1).Value = "even_MAG-GUT91328.fa"
$ws.Cells.Item(9,2).Value = 0.5347617402372367
$ws.Cells.Item(9,3).Value = 0.001064837127594029
$ws.Cells.Item(9,4).Value = 0.0003099261727542438
$ws.Cells.Item(9,5).Value = 0.05588100135301438
$ws.Cells.Item(9,6).Value = [double]"1.127164823491758e-13"
$ws.Cells.Item(9,7).Value = 0.03681116754575632
$ws.Cells.Item(9,8).Value = 0.002365986204945783
$ws.Cells.Item(9,9).Value = [double]"1.127164823491758e-13"
$ws.Cells.Item(9,10).Value = [double]"1.127164823491758e-13"
$ws.Cells.Item(9,11).Value = 0.1390495590006426
$ws.Cells.Item(9,12).Value = 0.1970557140393906
$ws.Cells.Item(9,13).Value = 0.005752349389228893
$ws.Cells.Item(9,14).Value = 0.02694771892909861
$ws.Cells.Item(9,15).Value = 0.5347617402372367
$ws.Cells.Item(9,16).Value = "g__Anaerococcus"
$ws.Cells.Item(9,17).Value = "g__Anaerococcus(reject)"

# Row 10: even_MAG-GUT91672.fa
$ws.Cells.Item(10,1).Value = "even_MAG-GUT91672.fa"
$ws.Cells.Item(10,2).Value = 0.00146909967817452
$ws.Cells.Item(10,3).Value = [double]"2.609752242150503e-05"
$ws.Cells.Item(10,4).Value = [double]"5.21845813186488e-07"
$ws.Cells.Item(10,5).Value = 0.002454492137544975
$ws.Cells.Item(10,6).Value = [double]"2.722059842679771e-14"
$ws.Cells.Item(10,7).Value = 0.006635920093821569
$ws.Cells.Item(10,8).Value = 0.001214640584467936
$ws.Cells.Item(10,9).Value = [double]"2.722059842679773e-14"
$ws.Cells.Item(10,10).Value = [double]"2.72205984267977e-14"
$ws.Cells.Item(10,11).Value = 0.02911580679341148
$ws.Cells.Item(10,12).Value = 0.9462993410575473
$ws.Cells.Item(10,13).Value = [double]"4.243539665185686e-05"
$ws.Cells.Item(10,14).Value = 0.01274164489006397
$ws.Cells.Item(10,15).Value = 0.9462993410575473
$ws.Cells.Item(10,16).Value = "g__Peptoniphilus_A"
$ws.Cells.Item(10,17).Value = "g__Peptoniphilus_A"

# Row 11: even_MAG-GUT91675.fa
$ws.Cells.Item(11,1).Value = "even_MAG-GUT91675.fa"
$ws.Cells.Item(11,2).Value = 0.002413487586914751
$ws.Cells.Item(11,3).Value = [double]"1.172355136776727e-05"
$ws.Cells.Item(11,4).Value = [double]"2.593640260629962e-07"
$ws.Cells.Item(11,5).Value = 0.004327063184864447
$ws.Cells.Item(11,6).Value = [double]"4.211629435226421e-14"
$ws.Cells.Item(11,7).Value = 0.009785633746628346
$ws.Cells.Item(11,8).Value = 0.001944334709166353
$ws.Cells.Item(11,9).Value = [double]"4.211629435226423e-14"
$ws.Cells.Item(11,10).Value = [double]"4.211629435226421e-14"
$ws.Cells.Item(11,11).Value = 0.05756883998352685
$ws.Cells.Item(11,12).Value = 0.9075829925404758
$ws.Cells.Item(11,13).Value = [double]"6.310910974378769e-05"
$ws.Cells.Item(11,14).Value = 0.01630255622315954
$ws.Cells.Item(11,15).Value = 0.9075829925404758
$ws.Cells.Item(11,16).Value = "g__Peptoniphilus_A"
$ws.Cells.Item(11,17).Value = "g__Peptoniphilus_A"
